$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D4").Value = 1797
$ws.Range("D5").Value = "19.8 (21.6)"
$ws.Range("D28").Value = "72 (4.1)"
$ws.Range("D29").Value = "628 (35.8)"
$ws.Range("D30").Value = "581 (33.1)"
$ws.Range("D31").Value = "225 (12.8)"
$ws.Range("D35").Value = "15 (0.9)"
$ws.Range("D37").ClearContents()
$ws.Range("D38").Value = "194 (11.0)"
